$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C holds the "Förändrad" (changed) date, stored as an Excel
# serial date number. Find the last used row in column A (the
# "Beteckning" key column) so the update covers the whole data range.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row()
if ($lastRow -lt 2) {
    $lastRow = 2
}

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 3).Value = 45204
}
